$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.620801766869107
$ws.Range("C2").Value = 0.2311657902019135
$ws.Range("D2").Value = 0.047776838161397
$ws.Range("F2").Value = 1.01695826612783
$ws.Range("G2").Value = 0.002467076477580325
$ws.Range("K2").Value = 0.2935921447024441
$ws.Range("L2").Value = 0.3037098705232211
$ws.Range("O2").Value = 3.642715521312738
$ws.Range("B3").Value = 0.5768644779398358
$ws.Range("C3").Value = 0.2327238799209788
$ws.Range("D3").Value = 0.04606343661111367
$ws.Range("F3").Value = 1.016913411053267
$ws.Range("G3").Value = 0.002469692833632317
$ws.Range("K3").Value = 0.2586924008109008
$ws.Range("L3").Value = 0.2926612632286236
$ws.Range("O3").Value = 3.657195331644687
$ws.Range("B4").Value = 0.5500901395506048
$ws.Range("C4").Value = 0.2337315644787061
$ws.Range("D4").Value = 0.04499954316590049
$ws.Range("F4").Value = 1.017462858044148
$ws.Range("G4").Value = 0.002471385406385599
$ws.Range("K4").Value = 0.2372526210518231
$ws.Range("L4").Value = 0.2860523327749576
$ws.Range("O4").Value = 3.668019184729502
$ws.Range("B5").Value = 0.5392311026547816
$ws.Range("C5").Value = 0.2341550536282604
$ws.Range("D5").Value = 0.04456303895602787
$ws.Range("F5").Value = 1.017831844069995
$ws.Range("G5").Value = 0.002472096863621177
$ws.Range("K5").Value = 0.228513371881462
$ws.Range("L5").Value = 0.2834031636323147
$ws.Range("O5").Value = 3.672916010886524
$ws.Range("B6").Value = 0.5374311093463575
$ws.Range("C6").Value = 0.2342261504188592
$ws.Range("D6").Value = 0.04449037972236169
$ws.Range("F6").Value = 1.01790187632772
$ws.Range("G6").Value = 0.002472216314087658
$ws.Range("K6").Value = 0.2270620967027526
$ws.Range("L6").Value = 0.2829659309021224
$ws.Range("O6").Value = 3.673758476585249
$ws.Range("B7").Value = 0.5499434805666397
$ws.Range("C7").Value = 0.2337372237458748
$ws.Range("D7").Value = 0.0449936682689227
$ws.Range("F7").Value = 1.017467246898036
$ws.Range("G7").Value = 0.002471394913372837
$ws.Range("K7").Value = 0.2371347694570574
$ws.Range("L7").Value = 0.2860164268643786
$ws.Range("O7").Value = 3.668083257308524
$ws.Range("B8").Value = 0.6056103281779031
$ws.Range("C8").Value = 0.231692445672131
$ws.Range("D8").Value = 0.0471885315891214
$ws.Range("F8").Value = 1.01682303953708
$ws.Range("G8").Value = 0.002467960762054227
$ws.Range("K8").Value = 0.2815612955741926
$ws.Range("L8").Value = 0.2998640110584887
$ws.Range("O8").Value = 3.647306861604648
$ws.Range("B9").Value = 0.7163660987903029
$ws.Range("C9").Value = 0.228086351999286
$ws.Range("D9").Value = 0.05139774925337548
$ws.Range("F9").Value = 1.020139572358545
$ws.Range("G9").Value = 0.002461906771343956
$ws.Range("K9").Value = 0.3685772851766558
$ws.Range("L9").Value = 0.3284078317835224
$ws.Range("O9").Value = 3.621910706266078
$ws.Range("B10").Value = 0.7986913851204349
$ws.Range("C10").Value = 0.2256815797370884
$ws.Range("D10").Value = 0.05443158504165524
$ws.Range("F10").Value = 1.025372541160621
$ws.Range("G10").Value = 0.002457869540390498
$ws.Range("K10").Value = 0.4324299632121154
$ws.Range("L10").Value = 0.3502289617875789
$ws.Range("O10").Value = 3.61262188013302
$ws.Range("B11").Value = 0.8363469254496749
$ws.Range("C11").Value = 0.2246403900176759
$ws.Range("D11").Value = 0.05579885864963785
$ws.Range("F11").Value = 1.028361401267958
$ws.Range("G11").Value = 0.002456121185681726
$ws.Range("K11").Value = 0.4614586389926103
$ws.Range("L11").Value = 0.3603414036310966
$ws.Range("O11").Value = 3.610434050262654
$ws.Range("B12").Value = 0.8506351574025643
$ws.Range("C12").Value = 0.2242536822173431
$ws.Range("D12").Value = 0.05631474496534139
$ws.Range("F12").Value = 1.029580732390272
$ws.Range("G12").Value = 0.002455471746027459
$ws.Range("K12").Value = 0.4724480597487855
$ws.Range("L12").Value = 0.364197467199574
$ws.Range("O12").Value = 3.609898802144954
$ws.Range("B13").Value = 0.8475566559087611
$ws.Range("C13").Value = 0.2243366303438066
$ws.Range("D13").Value = 0.05620372320677802
$ws.Range("F13").Value = 1.029314234851796
$ws.Range("G13").Value = 0.002455611053970663
$ws.Range("K13").Value = 0.4700814373647688
$ws.Range("L13").Value = 0.3633658078662734
$ws.Range("O13").Value = 3.610001032244526
$ws.Range("B14").Value = 0.8375218513941149
$ws.Range("C14").Value = 0.2246084237575232
$ws.Range("D14").Value = 0.05584133855808915
$ws.Range("F14").Value = 1.028459962259532
$ws.Range("G14").Value = 0.002456067503205372
$ws.Range("K14").Value = 0.462362810341034
$ws.Range("L14").Value = 0.3606581093463888
$ws.Range("O14").Value = 3.610384136967298
$ws.Range("B15").Value = 0.8313789882302558
$ws.Range("C15").Value = 0.2247758901485426
$ws.Range("D15").Value = 0.05561912315783246
$ws.Range("F15").Value = 1.027948093606255
$ws.Range("G15").Value = 0.002456348733979873
$ws.Range("K15").Value = 0.4576345091342944
$ws.Range("L15").Value = 0.3590030421921995
$ws.Range("O15").Value = 3.610656993990517
$ws.Range("B16").Value = 0.7962345850429529
$ws.Range("C16").Value = 0.2257506839437191
$ws.Range("D16").Value = 0.05434197010532671
$ws.Range("F16").Value = 1.025189456719474
$ws.Range("G16").Value = 0.002457985569261346
$ws.Range("K16").Value = 0.4305324624536127
$ws.Range("L16").Value = 0.3495718280883153
$ws.Range("O16").Value = 3.612805882290274
$ws.Range("B17").Value = 0.7747267766842185
$ws.Range("C17").Value = 0.2263621870308477
$ws.Range("D17").Value = 0.05355517321023484
$ws.Range("F17").Value = 1.023652961369535
$ws.Range("G17").Value = 0.002459012263600355
$ws.Range("K17").Value = 0.4139012183718478
$ws.Range("L17").Value = 0.3438336753974625
$ws.Range("O17").Value = 3.614646214842566
$ws.Range("B18").Value = 0.7623754208694606
$ws.Range("C18").Value = 0.2267188738552282
$ws.Range("D18").Value = 0.05310142198424472
$ws.Range("F18").Value = 1.022826469307134
$ws.Range("G18").Value = 0.002459611096268315
$ws.Range("K18").Value = 0.4043336799740587
$ws.Range("L18").Value = 0.340550741444801
$ws.Range("O18").Value = 3.615896510915292
$ws.Range("B19").Value = 0.7581968096094158
$ws.Range("C19").Value = 0.2268404954175995
$ws.Range("D19").Value = 0.0529475831713242
$ws.Range("F19").Value = 1.022556466613381
$ws.Range("G19").Value = 0.002459815279021504
$ws.Range("K19").Value = 0.4010940010732611
$ws.Range("L19").Value = 0.3394422029416404
$ws.Range("O19").Value = 3.616352772417542
$ws.Range("B20").Value = 0.7770143198191022
$ws.Range("C20").Value = 0.2262965776077515
$ws.Range("D20").Value = 0.05363905412877301
$ws.Range("F20").Value = 1.023810597537377
$ws.Range("G20").Value = 0.002458902111107553
$ws.Range("K20").Value = 0.4156718208014638
$ws.Range("L20").Value = 0.3444427006171935
$ws.Range("O20").Value = 3.614430457645511
$ws.Range("B21").Value = 0.8404685375355712
$ws.Range("C21").Value = 0.2245283861858418
$ws.Range("D21").Value = 0.05594783067847686
$ws.Range("F21").Value = 1.028708507441138
$ws.Range("G21").Value = 0.00245593309080462
$ws.Range("K21").Value = 0.4646300464455635
$ws.Range("L21").Value = 0.36145270162028
$ws.Range("O21").Value = 3.610263649778688
$ws.Range("B22").Value = 0.8821076728160051
$ws.Range("C22").Value = 0.2234168758702602
$ws.Range("D22").Value = 0.05744583523625835
$ws.Range("F22").Value = 1.03241969205483
$ws.Range("G22").Value = 0.002454066222270949
$ws.Range("K22").Value = 0.4966087322476653
$ws.Range("L22").Value = 0.3727253415656691
$ws.Range("O22").Value = 3.609249700567887
$ws.Range("B23").Value = 0.8598689070355476
$ws.Range("C23").Value = 0.2240060799957035
$ws.Range("D23").Value = 0.05664732892640245
$ws.Range("F23").Value = 1.030392276337139
$ws.Range("G23").Value = 0.002455055893886595
$ws.Range("K23").Value = 0.4795429448531081
$ws.Range("L23").Value = 0.3666946924907535
$ws.Range("O23").Value = 3.609634392177753
$ws.Range("B24").Value = 0.7759800786303401
$ws.Range("C24").Value = 0.2263262236591892
$ws.Range("D24").Value = 0.05360113595992289
$ws.Range("F24").Value = 1.023739153112885
$ws.Range("G24").Value = 0.002458951884315774
$ws.Range("K24").Value = 0.4148713501430734
$ws.Range("L24").Value = 0.3441673104572942
$ws.Range("O24").Value = 3.614527402589601
$ws.Range("B25").Value = 0.6862350576245433
$ws.Range("C25").Value = 0.229018831036738
$ws.Range("D25").Value = 0.05026928983853907
$ws.Range("F25").Value = 1.018751609995078
$ws.Range("G25").Value = 0.002463472127275014
$ws.Range("K25").Value = 0.3450498170722369
$ws.Range("L25").Value = 0.3205369474227524
$ws.Range("O25").Value = 3.627136639535479
